$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.107.89'
$ws.Range("E2").Value = '  +5.67%  '

$ws.Range("D3").Value = '1.913.41'
$ws.Range("E3").Value = '  +2.30%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.43%  '

$ws.Range("D5").Value = '330.28'
$ws.Range("E5").Value = '  +5.04%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.46%  '

$ws.Range("D7").Value = '0.5190'
$ws.Range("E7").Value = '  +2.41%  '

$ws.Range("D8").Value = '0.4060'
$ws.Range("E8").Value = '  +4.02%  '

$ws.Range("D9").Value = '0.08490'
$ws.Range("E9").Value = '  +1.90%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.125'
$ws.Range("E10").Value = '  +1.89%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '42.76'
$ws.Range("E11").Value = '  +0.63%  '

$ws.Range("D12").Value = '23.07'
$ws.Range("E12").Value = '  +13.69%  '

$ws.Range("D13").Value = '6.438'
$ws.Range("E13").Value = '  +4.03%  '

$ws.Range("D14").Value = '1.901.12'
$ws.Range("E14").Value = '  +1.66%  '

$ws.Range("E15").Value = '  +1.85%  '

$ws.Range("E16").Value = '  -0.46%  '

$ws.Range("D17").Value = '95.01'
$ws.Range("E17").Value = '  +2.66%  '

$ws.Range("D19").Value = '0.06694'
$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("D20").Value = '18.39'
$ws.Range("E20").Value = '  +4.46%  '

$ws.Range("E21").Value = '  -0.47%  '

$ws.Range("D22").Value = '6.007'
$ws.Range("E22").Value = '  +1.59%  '

$ws.Range("D23").Value = '30.137.37'
$ws.Range("E23").Value = '  +5.60%  '

$ws.Range("D24").Value = '11.33'
$ws.Range("E24").Value = '  +2.56%  '

$ws.Range("D25").Value = '2.228'
$ws.Range("E25").Value = '  +1.79%  '

$ws.Range("D26").Value = '2.134.31'
$ws.Range("E26").Value = '  +2.47%  '

$ws.Range("D27").Value = '161.64'
$ws.Range("E27").Value = '  +2.43%  '

$ws.Range("D28").Value = '21.08'
$ws.Range("E28").Value = '  +2.65%  '

$ws.Range("D29").Value = '2.407'
$ws.Range("E29").Value = '  -0.39%  '

$ws.Range("D30").Value = '128.38'
$ws.Range("E30").Value = '  +1.99%  '

$ws.Range("D31").Value = '1.101'
$ws.Range("E31").Value = '  +6.00%  '

$ws.Range("D32").Value = '0.1066'
$ws.Range("E32").Value = '  +3.16%  '

$ws.Range("D33").Value = '5.992'
$ws.Range("E33").Value = '  +3.95%  '

$ws.Range("D34").Value = '3.621'
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("D35").Value = '0.02492'
$ws.Range("E35").Value = '  +2.04%  '

$ws.Range("D36").Value = '0.06563'
$ws.Range("E36").Value = '  +0.29%  '

$ws.Range("D37").Value = '0.2210'
$ws.Range("E37").Value = '  +2.41%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '5.162'
$ws.Range("E38").Value = '  +2.76%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = '1.224'
$ws.Range("E39").Value = '  +3.48%  '

$ws.Range("D40").Value = '11.89'
$ws.Range("E40").Value = '  +7.41%  '

$ws.Range("D41").Value = '8.805'
$ws.Range("E41").Value = '  -1.81%  '

$ws.Range("D42").Value = '0.6511'
$ws.Range("E42").Value = '  +2.48%  '

$ws.Range("D43").Value = '1.239'
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").Value = '0.6128'
$ws.Range("E44").Value = '  +2.49%  '

$ws.Range("D45").Value = '13.32'
$ws.Range("E45").Value = '  +2.02%  '

$ws.Range("D46").Value = '3.747'
$ws.Range("E46").Value = '  +1.94%  '

$ws.Range("D47").Value = '2.078'
$ws.Range("E47").Value = '  +4.01%  '

$ws.Range("D48").Value = '1.240'
$ws.Range("E48").Value = '  +2.63%  '

$ws.Range("D49").Value = '123.86'
$ws.Range("E49").Value = '  +1.43%  '

$ws.Range("E50").Value = '  +0.45%  '

$ws.Range("D51").Value = '79.32'
$ws.Range("E51").Value = '  +4.32%  '
